$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "1.004", "151.00") are preserved exactly as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.887.71'
$ws.Range("E2").Value = '  +2.50%  '

$ws.Range("D3").Value = '1.715.04'
$ws.Range("E3").Value = '  +2.56%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '311.46'
$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.41%  '

$ws.Range("D7").Value = '0.3774'
$ws.Range("E7").Value = '  +1.21%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3487'
$ws.Range("E8").Value = '  +1.83%  '

$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '49.48'
$ws.Range("E9").Value = '  +3.99%  '

$ws.Range("D10").Value = '1.194'
$ws.Range("E10").Value = '  +1.07%  '

$ws.Range("D11").Value = '0.07485'
$ws.Range("E11").Value = '  +2.81%  '

$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.20%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '20.91'
$ws.Range("E13").Value = '  +2.32%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '6.290'
$ws.Range("E14").Value = '  +3.32%  '

$ws.Range("D15").Value = '6.991'
$ws.Range("E15").Value = '  +3.63%  '

$ws.Range("D16").Value = '1.720.38'
$ws.Range("E16").Value = '  +2.63%  '

$ws.Range("D17").Value = '0.00001126'
$ws.Range("E17").Value = '  +1.79%  '

$ws.Range("E18").Value = '  +0.37%  '

$ws.Range("D19").Value = '0.06732'
$ws.Range("E19").Value = '  +0.31%  '

$ws.Range("D20").Value = '84.22'
$ws.Range("E20").Value = '  +3.40%  '

$ws.Range("D21").Value = '17.24'
$ws.Range("E21").Value = '  +5.09%  '

$ws.Range("D22").Value = '6.395'
$ws.Range("E22").Value = '  +4.27%  '

$ws.Range("D23").Value = '13.04'
$ws.Range("E23").Value = '  +8.37%  '

$ws.Range("D24").Value = '24.828.95'
$ws.Range("E24").Value = '  +2.49%  '

$ws.Range("D25").Value = '2.436'
$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("D26").Value = '2.792'
$ws.Range("E26").Value = '  +5.39%  '

$ws.Range("D27").Value = '20.53'
$ws.Range("E27").Value = '  +5.53%  '

$ws.Range("D28").Value = '151.00'
$ws.Range("E28").Value = '  -0.49%  '

$ws.Range("D29").Value = '131.79'
$ws.Range("E29").Value = '  +3.88%  '

$ws.Range("D30").Value = '1.909.11'
$ws.Range("E30").Value = '  +2.62%  '

$ws.Range("D31").Value = '1.178'
$ws.Range("E31").Value = '  +19.31%  '

$ws.Range("D32").Value = '6.806'
$ws.Range("E32").Value = '  +6.80%  '

$ws.Range("D33").Value = '4.243'
$ws.Range("E33").Value = '  +4.55%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '13.79'
$ws.Range("E34").Value = '  +9.98%  '

$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.08845'
$ws.Range("E35").Value = '  +4.64%  '

$ws.Range("D36").Value = '1.770'
$ws.Range("E36").Value = '  +1.08%  '

$ws.Range("D37").Value = '5.610'
$ws.Range("E37").Value = '  +4.76%  '

$ws.Range("D38").Value = '0.06554'
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("D39").Value = '0.02409'
$ws.Range("E39").Value = '  +2.94%  '

$ws.Range("D40").Value = '8.983'
$ws.Range("E40").Value = '  +2.06%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2202'
$ws.Range("E41").Value = '  +4.59%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.283'
$ws.Range("E42").Value = '  -0.06%  '

$ws.Range("D43").Value = '0.6439'
$ws.Range("E43").Value = '  +4.69%  '

$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.30%  '

$ws.Range("D45").Value = '13.92'
$ws.Range("E45").Value = '  +4.61%  '

$ws.Range("D46").Value = '0.6147'
$ws.Range("E46").Value = '  +3.46%  '

$ws.Range("D47").Value = '3.831'
$ws.Range("E47").Value = '  +0.78%  '

$ws.Range("D48").Value = '2.143'
$ws.Range("E48").Value = '  +6.30%  '

$ws.Range("D49").Value = '130.14'
$ws.Range("E49").Value = '  +2.00%  '

$ws.Range("D50").Value = '0.07267'
$ws.Range("E50").Value = '  +1.54%  '

$ws.Range("D51").Value = '79.82'
$ws.Range("E51").Value = '  +4.24%  '
